# Update attendee counts ("想去人数") that changed between data refreshes.
# Sheet "展览" (worksheet 1): row 3 -> 209, row 4 -> 827, row 5 -> 72
# Sheet "全部类型" (worksheet 4): row 4 -> 209, row 5 -> 827, row 6 -> 72

$wb = $excel.ActiveWorkbook

$wsExhibition = $wb.Worksheets.Item("展览")
$wsExhibition.Range("F3").Value = 209
$wsExhibition.Range("F4").Value = 827
$wsExhibition.Range("F5").Value = 72

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value = 209
$wsAll.Range("F5").Value = 827
$wsAll.Range("F6").Value = 72
